$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1326
$ws.Range("I19").Value = 999.5714
$ws.Range("J19").Value = 1897.25
$ws.Range("K19").Value = 999.5714
$ws.Range("L19").Value = 1897.25
$ws.Range("M19").Value = -824.5714
$ws.Range("N19").Value = -2247.25
$ws.Range("H70").Value = 2805.5
$ws.Range("I70").Value = 1500
$ws.Range("J70").Value = 3066.6
$ws.Range("K70").Value = 4500
$ws.Range("L70").Value = 9199.799999999999
$ws.Range("M70").Value = -4230
$ws.Range("N70").Value = -9739.799999999999
$ws.Range("H73").Value = 2805.5
$ws.Range("I73").Value = 1500
$ws.Range("J73").Value = 3066.6
$ws.Range("K73").Value = 4500
$ws.Range("L73").Value = 9199.799999999999
$ws.Range("M73").Value = -3564
$ws.Range("N73").Value = -11071.8
$ws.Range("H98").Value = 449053.03
$ws.Range("I98").Value = 746466.9399999999
$ws.Range("K98").Value = 746466.9399999999
$ws.Range("M98").Value = -744968.9399999999
$ws.Range("H113").Value = 86444.836
$ws.Range("I113").Value = 114204.22
$ws.Range("J113").Value = 3166.6667
$ws.Range("K113").Value = 114204.22
$ws.Range("L113").Value = 3166.6667
$ws.Range("M113").Value = -110950.22
$ws.Range("N113").Value = -9674.6667
$ws.Range("H122").Value = 449053.03
$ws.Range("I122").Value = 746466.9399999999
$ws.Range("K122").Value = 2239400.82
$ws.Range("M122").Value = -2236950.82
$ws.Range("H132").Value = 445628.72
$ws.Range("I132").Value = 487103.9
$ws.Range("J132").Value = 100002
$ws.Range("K132").Value = 1461311.7
$ws.Range("L132").Value = 300006
$ws.Range("M132").Value = -1458781.7
$ws.Range("N132").Value = -305066
$ws.Range("H135").Value = 1209.3334
$ws.Range("I135").Value = 1102.0714
$ws.Range("J135").Value = 1960.1666
$ws.Range("K135").Value = 9918.642600000001
$ws.Range("L135").Value = 17641.4994
$ws.Range("M135").Value = -7383.642600000001
$ws.Range("N135").Value = -22711.4994
$ws.Range("H137").Value = 24391486
$ws.Range("I137").Value = 33334266
$ws.Range("K137").Value = 100002798
$ws.Range("M137").Value = -100000248

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 64957.125
$ws.Range("I2").Value = 92783.09
$ws.Range("K2").Value = 92783.09
$ws.Range("M2").Value = -92670.09
$ws.Range("H32").Value = 21065.715
$ws.Range("I32").Value = 2955.9622
$ws.Range("K32").Value = 2955.9622
$ws.Range("M32").Value = -2668.9622
$ws.Range("H45").Value = 1028
$ws.Range("I45").Value = 1032
$ws.Range("J45").Value = 1000
$ws.Range("K45").Value = 1032
$ws.Range("L45").Value = 1000
$ws.Range("M45").Value = -655
$ws.Range("N45").Value = -1754
$ws.Range("H61").Value = 2200.4048
$ws.Range("I61").Value = 1550.8387
$ws.Range("J61").Value = 4031
$ws.Range("K61").Value = 1550.8387
$ws.Range("L61").Value = 4031
$ws.Range("M61").Value = -1338.8387
$ws.Range("N61").Value = -4455
$ws.Range("H74").Value = 4215.244
$ws.Range("I74").Value = 1027.2059
$ws.Range("J74").Value = 19700
$ws.Range("K74").Value = 1027.2059
$ws.Range("L74").Value = 19700
$ws.Range("M74").Value = -153.2058999999999
$ws.Range("N74").Value = -21448
$ws.Range("H77").Value = 4215.244
$ws.Range("I77").Value = 1027.2059
$ws.Range("J77").Value = 19700
$ws.Range("K77").Value = 5136.0295
$ws.Range("L77").Value = 98500
$ws.Range("M77").Value = -768.0294999999996
$ws.Range("N77").Value = -107236
$ws.Range("H116").Value = 64957.125
$ws.Range("I116").Value = 92783.09
$ws.Range("K116").Value = 92783.09
$ws.Range("M116").Value = -90489.09
$ws.Range("H132").Value = 2623.5112
$ws.Range("I132").Value = 2152.8125
$ws.Range("J132").Value = 3782.1538
$ws.Range("K132").Value = 6458.4375
$ws.Range("L132").Value = 11346.4614
$ws.Range("M132").Value = -3928.4375
$ws.Range("N132").Value = -16406.4614
$ws.Range("H136").Value = 2200.4048
$ws.Range("I136").Value = 1550.8387
$ws.Range("J136").Value = 4031
$ws.Range("K136").Value = 4652.5161
$ws.Range("L136").Value = 12093
$ws.Range("M136").Value = -2102.5161
$ws.Range("N136").Value = -17193

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 64957.125
$ws.Range("I3").Value = 92783.09
$ws.Range("K3").Value = 92783.09
$ws.Range("M3").Value = -92669.09
$ws.Range("H86").Value = 11863.818
$ws.Range("I86").Value = 3799
$ws.Range("J86").Value = 18584.5
$ws.Range("K86").Value = 3799
$ws.Range("L86").Value = 18584.5
$ws.Range("M86").Value = -2676
$ws.Range("N86").Value = -20830.5
$ws.Range("H89").Value = 11863.818
$ws.Range("I89").Value = 3799
$ws.Range("J89").Value = 18584.5
$ws.Range("K89").Value = 18995
$ws.Range("L89").Value = 92922.5
$ws.Range("M89").Value = -13379
$ws.Range("N89").Value = -104154.5
$ws.Range("H107").Value = 1112.4445
$ws.Range("I107").Value = 1051.5
$ws.Range("J107").Value = 1600
$ws.Range("K107").Value = 1051.5
$ws.Range("L107").Value = 1600
$ws.Range("M107").Value = 868.5
$ws.Range("N107").Value = -5440
$ws.Range("H134").Value = 2585.6606
$ws.Range("I134").Value = 1666.3334
$ws.Range("J134").Value = 4694.706
$ws.Range("K134").Value = 4999.0002
$ws.Range("L134").Value = 14084.118
$ws.Range("M134").Value = -2464.0002
$ws.Range("N134").Value = -19154.118

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 2584.6
$ws.Range("I5").Value = 238.33333
$ws.Range("J5").Value = 6104
$ws.Range("K5").Value = 238.33333
$ws.Range("L5").Value = 6104
$ws.Range("M5").Value = -126.33333
$ws.Range("N5").Value = -6328
$ws.Range("H31").Value = 1688.186
$ws.Range("I31").Value = 910.3214
$ws.Range("J31").Value = 3140.2
$ws.Range("K31").Value = 910.3214
$ws.Range("L31").Value = 3140.2
$ws.Range("M31").Value = -615.3214
$ws.Range("N31").Value = -3730.2
$ws.Range("H34").Value = 1688.186
$ws.Range("I34").Value = 910.3214
$ws.Range("J34").Value = 3140.2
$ws.Range("K34").Value = 910.3214
$ws.Range("L34").Value = 3140.2
$ws.Range("M34").Value = -708.3214
$ws.Range("N34").Value = -3544.2
$ws.Range("H132").Value = 1927.5555
$ws.Range("I132").Value = 1470.65
$ws.Range("J132").Value = 3233
$ws.Range("K132").Value = 4411.950000000001
$ws.Range("L132").Value = 9699
$ws.Range("M132").Value = -1881.950000000001
$ws.Range("N132").Value = -14759
$ws.Range("H134").Value = 2220.44
$ws.Range("I134").Value = 1445.4054
$ws.Range("J134").Value = 4426.3076
$ws.Range("K134").Value = 4336.216200000001
$ws.Range("L134").Value = 13278.9228
$ws.Range("M134").Value = -1801.216200000001
$ws.Range("N134").Value = -18348.9228

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 445
$ws.Range("I107").Value = 443
$ws.Range("J107").Value = 455
$ws.Range("K107").Value = 1329
$ws.Range("L107").Value = 1365
$ws.Range("M107").Value = 591
$ws.Range("N107").Value = -5205
$ws.Range("H140").Value = 3829.5
$ws.Range("I140").Value = 4692.407
$ws.Range("K140").Value = 14077.221
$ws.Range("M140").Value = -8897.221000000001

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 7000
$ws.Range("I43").Value = 6500
$ws.Range("J43").Value = 8000
$ws.Range("K43").Value = 6500
$ws.Range("L43").Value = 8000
$ws.Range("M43").Value = -6349
$ws.Range("N43").Value = -8302
$ws.Range("H46").Value = 29000
$ws.Range("I46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("M46").ClearContents()
$ws.Range("H102").Value = 2651.2942
$ws.Range("I102").Value = 2936
$ws.Range("J102").Value = 2244.5715
$ws.Range("K102").Value = 2936
$ws.Range("L102").Value = 2244.5715
$ws.Range("M102").Value = -1314
$ws.Range("N102").Value = -5488.5715
$ws.Range("H126").Value = 2144.8076
$ws.Range("I126").Value = 1583.125
$ws.Range("J126").Value = 2394.4443
$ws.Range("K126").Value = 4749.375
$ws.Range("L126").Value = 7183.3329
$ws.Range("M126").Value = -2279.375
$ws.Range("N126").Value = -12123.3329
$ws.Range("H132").Value = 2467.8386
$ws.Range("I132").Value = 2324.5
$ws.Range("J132").Value = 2818.2222
$ws.Range("K132").Value = 6973.5
$ws.Range("L132").Value = 8454.6666
$ws.Range("M132").Value = -4443.5
$ws.Range("N132").Value = -13514.6666

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1429.3334
$ws.Range("I93").Value = 543.2
$ws.Range("J93").Value = 2537
$ws.Range("K93").Value = 543.2
$ws.Range("L93").Value = 2537
$ws.Range("M93").Value = 704.8
$ws.Range("N93").Value = -5033
$ws.Range("H132").Value = 2582.7368
$ws.Range("I132").Value = 2034.7368
$ws.Range("J132").Value = 3678.7368
$ws.Range("K132").Value = 6104.2104
$ws.Range("L132").Value = 11036.2104
$ws.Range("M132").Value = -3574.2104
$ws.Range("N132").Value = -16096.2104
$ws.Range("H136").Value = 3959.3262
$ws.Range("I136").Value = 2411.1943
$ws.Range("J136").Value = 9532.6
$ws.Range("K136").Value = 7233.5829
$ws.Range("L136").Value = 28597.8
$ws.Range("M136").Value = -4683.5829
$ws.Range("N136").Value = -33697.8

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H19").Value = 50000
$ws.Range("J19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("N19").ClearContents()
$ws.Range("H122").Value = 93407.45
$ws.Range("I122").Value = 201896.8
$ws.Range("J122").Value = 2999.6667
$ws.Range("K122").Value = 605690.3999999999
$ws.Range("L122").Value = 8999.000100000001
$ws.Range("M122").Value = -603240.3999999999
$ws.Range("N122").Value = -13899.0001
$ws.Range("H126").Value = 201421
$ws.Range("I126").Value = 251276.25
$ws.Range("K126").Value = 753828.75
$ws.Range("M126").Value = -751358.75
$ws.Range("H128").Value = 73973.44
$ws.Range("J128").Value = 73973.44
$ws.Range("L128").Value = 73973.44
$ws.Range("N128").Value = -83933.44
$ws.Range("H132").Value = 7938213.5
$ws.Range("I132").Value = 12196476
$ws.Range("K132").Value = 36589428
$ws.Range("M132").Value = -36586898
$ws.Range("H136").Value = 5481768.5
$ws.Range("I136").Value = 6430260
$ws.Range("J136").Value = 1593.7778
$ws.Range("K136").Value = 19290780
$ws.Range("L136").Value = 4781.3334
$ws.Range("M136").Value = -19288230
$ws.Range("N136").Value = -9881.3334
